# Update countries & provincias Spain
# - Reorder two pairs of countries alphabetically (names swap places,
#   taking their row's statistics with them):
#     row 181 <-> row 182  : "Islas Turcas y Caicos" / "San Martin (Parte Holandesa)"
#     row 213 <-> row 214  : "Islas Malvinas" / "Montserrat"
# - Refresh the COVID-19 case counters for several countries.
# - Bump the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swap: Islas Turcas y Caicos <-> San Martin (Parte Holandesa) ---
$ws.Range("A181").Value = "San Martin (Parte Holandesa)"
$ws.Range("A182").Value = "Islas Turcas y Caicos"

# --- Country name swap: Islas Malvinas <-> Montserrat ---
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Muertes) ---

# Australia (row 72)
$ws.Range("B72").Value = 22356
$ws.Range("C72").Value = 229
$ws.Range("E72").Value = 9216

# Jamaica (row 154)
$ws.Range("B154").Value = 1065
$ws.Range("C154").Value = 18
$ws.Range("E154").Value = 298

# Mongolia (row 177)
$ws.Range("B177").Value = 297
$ws.Range("C177").Value = 4
$ws.Range("E177").Value = 28

# Camboya (row 180)
$ws.Range("B180").Value = 272
$ws.Range("C180").Value = 4
$ws.Range("D180").Value = 223
$ws.Range("E180").Value = 49

# Row 181 -> now "San Martin (Parte Holandesa)"
$ws.Range("B181").Value = 248
$ws.Range("C181").Value = 29
$ws.Range("D181").Value = 102
$ws.Range("E181").Value = 129
$ws.Range("H181").Value = 17

# Row 182 -> now "Islas Turcas y Caicos"
$ws.Range("B182").Value = 241
$ws.Range("C182").Value = 17
$ws.Range("D182").Value = 46
$ws.Range("E182").Value = 193
$ws.Range("H182").Value = 2

# Row 213 -> now "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 -> now "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# --- Bump the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 06:01"
